$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "34.767.06"
$ws.Range("E2").Value = "  -1.74%  "

# Row 3
$ws.Range("D3").Value = "1.811.92"
$ws.Range("E3").Value = "  -1.36%  "

# Row 4
$ws.Range("E4").Value = "  +0.21%  "

# Row 5
$ws.Range("D5").Value = "'230.20"
$ws.Range("E5").Value = "  -0.09%  "

# Row 6
$ws.Range("D6").Value = "'0.608"
$ws.Range("E6").Value = "  -0.10%  "

# Row 7
$ws.Range("E7").Value = "  +0.21%  "

# Row 8
$ws.Range("D8").Value = "'39.74"
$ws.Range("E8").Value = "  -8.88%  "

# Row 9
$ws.Range("D9").Value = "'0.323"
$ws.Range("E9").Value = "  +5.36%  "

# Row 10
$ws.Range("E10").Value = "  -2.28%  "

# Row 11
$ws.Range("D11").Value = "'0.0994"
$ws.Range("E11").Value = "  -1.40%  "

# Row 12
$ws.Range("D12").Value = "2.073.65"
$ws.Range("E12").Value = "  -1.39%  "

# Row 13
$ws.Range("D13").Value = "'11.27"
$ws.Range("E13").Value = "  +0.37%  "

# Row 14
$ws.Range("D14").Value = "'0.666"
$ws.Range("E14").Value = "  -0.42%  "

# Row 15
$ws.Range("D15").Value = "1.807.80"
$ws.Range("E15").Value = "  -1.64%  "

# Row 16
$ws.Range("E16").Value = "  -1.29%  "

# Row 17
$ws.Range("D17").Value = "34.755.68"
$ws.Range("E17").Value = "  -1.76%  "

# Row 18
$ws.Range("D18").Value = "'69.51"
$ws.Range("E18").Value = "  -0.49%  "

# Row 19
$ws.Range("D19").Value = "0.0₃0784"
$ws.Range("E19").Value = "  -1.53%  "

# Row 20
$ws.Range("D20").Value = "'240.92"
$ws.Range("E20").Value = "  -1.37%  "

# Row 21
$ws.Range("E21").Value = "  -0.65%  "

# Row 22
$ws.Range("D22").Value = "'4.69"
$ws.Range("E22").Value = "  +0.70%  "

# Row 23
$ws.Range("E23").Value = "  +0.22%  "

# Row 24
$ws.Range("E24").Value = "  +2.24%  "

# Row 25
$ws.Range("D25").Value = "'171.63"
$ws.Range("E25").Value = "  +1.38%  "

# Row 26
$ws.Range("D26").Value = "'7.79"
$ws.Range("E26").Value = "  -1.12%  "

# Row 27
$ws.Range("D27").Value = "'17.29"
$ws.Range("E27").Value = "  -2.02%  "

# Row 28
$ws.Range("D28").Value = "'0.121"
$ws.Range("E28").Value = "  -0.23%  "

# Row 29
$ws.Range("E29").Value = "  +0.88%  "

# Row 30
$ws.Range("E30").Value = "  +0.18%  "

# Row 31
$ws.Range("D31").Value = "'4.08"
$ws.Range("E31").Value = "  +4.11%  "

# Row 32
$ws.Range("E32").Value = "  -0.29%  "

# Row 33
$ws.Range("E33").Value = "  -2.35%  "

# Row 34
$ws.Range("D34").Value = "'1.25"
$ws.Range("E34").Value = "  +15.28%  "

# Row 35
$ws.Range("E35").Value = "  -2.89%  "

# Row 36
$ws.Range("D36").Value = "'0.702"
$ws.Range("E36").Value = "  +2.99%  "

# Row 37
$ws.Range("D37").Value = "'92.03"
$ws.Range("E37").Value = "  -3.59%  "

# Row 38
$ws.Range("E38").Value = "  +4.43%  "

# Row 39
$ws.Range("D39").Value = "1.335.12"
$ws.Range("E39").Value = "  -0.46%  "

# Row 40
$ws.Range("D40").Value = "'0.0193"
$ws.Range("E40").Value = "  -0.69%  "

# Row 41
$ws.Range("B41").Value = "ARBITRUM"
$ws.Range("C41").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D41").Value = "'0.967"
$ws.Range("E41").Value = "  -3.56%  "

# Row 42
$ws.Range("B42").Value = "HuobiToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D42").Value = "'2.47"
$ws.Range("E42").Value = "  +0.27%  "

# Row 43
$ws.Range("D43").Value = "'14.31"
$ws.Range("E43").Value = "  -6.95%  "

# Row 44
$ws.Range("D44").Value = "'2.21"
$ws.Range("E44").Value = "  -8.13%  "

# Row 45
$ws.Range("D45").Value = "'2.69"
$ws.Range("E45").Value = "  -4.33%  "

# Row 46
$ws.Range("D46").Value = "'6.24"
$ws.Range("E46").Value = "  +0.66%  "

# Row 47
$ws.Range("E47").Value = "  -1.27%  "

# Row 48
$ws.Range("D48").Value = "2.000.47"
$ws.Range("E48").Value = "  -0.27%  "

# Row 49
$ws.Range("E49").Value = "  +0.16%  "

# Row 50
$ws.Range("D50").Value = "'0.0668"
$ws.Range("E50").Value = "  +6.78%  "

# Row 51
$ws.Range("D51").Value = "'98.04"
$ws.Range("E51").Value = "  -4.27%  "
